# Re-add the pruned rows of zeros (and their date stamp) to the simulation
# tracker sheet, plus the trailing marker cell that was dropped along with
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 35: new simulation entry -----------------------------------------
# Date stamp in column A. Copy the date formatting from an existing date
# cell (A3) instead of setting .NumberFormat directly, so the cell reuses
# the workbook's existing built-in date style instead of minting a new
# custom number format.
$ws.Range("A35").Value = 45203
$ws.Range("A3").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "rows of zeros" (really small leftover proportions) for the bins J:V.
$ws.Range("J35").Value = 0.07
$ws.Range("K35").Value = 0.04
$ws.Range("L35").Value = 0.07
$ws.Range("M35").Value = 0.04
$ws.Range("N35").Value = 0.04
$ws.Range("O35").Value = 0.05
$ws.Range("P35").Value = 0.03
$ws.Range("Q35").Value = 0.03
$ws.Range("R35").Value = 0.03
$ws.Range("S35").Value = 0.03
$ws.Range("T35").Value = 0.02
$ws.Range("U35").Value = 0.03
$ws.Range("V35").Value = 0.03

# --- Row 38: lone trailing marker cell -------------------------------------
$ws.Range("U38").Value = " "

# --- Selection / view bookkeeping ------------------------------------------
$ws.Activate()
$ws.Range("A36").Select()
